$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: G2 "Recorded By" list reordered ---
$ws.Range("G2").Value = "Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"

# --- Row 3: G3 "Recorded By" list reordered ---
$ws.Range("G3").Value = "hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 4: G4 "Recorded By" list reordered ---
$ws.Range("G4").Value = "gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 5: G5 "Recorded By" list reordered ---
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 6: G6 "Recorded By" list reordered, L6 Recorded Sessions 20 -> 21 ---
$ws.Range("G6").Value = "Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
$ws.Range("L6").Value = 21

# --- Row 7: G7 "Recorded By" list reordered ---
$ws.Range("G7").Value = "lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg"

# --- Row 8: L8 Pending Sessions 7 -> 6 ---
$ws.Range("L8").Value = 6

# --- Row 9: L9 Coverage % 69.0% -> 72.4% (force as literal text, not a numeric percent) ---
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "72.4%"

# --- Row 10: L10 Average Attendance % 26.5% -> 25.4% (force as literal text) ---
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "25.4%"

# --- Row 12: G12 "Recorded By" list reordered ---
$ws.Range("G12").Value = "Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"

# --- Row 15: G15 "Recorded By" list reordered, and the Class Statistics summary row (O/Q/R/S) updated to match L6/L8/L9/L10 ---
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("O15").Value = 21
$ws.Range("Q15").Value = 6
$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "72.4%"
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "25.4%"

# --- Row 21: PARASITOLOGY SGD/POS session was Pending, now Recorded ---
# Re-format A21:I21 to the normal (non-highlighted) look used elsewhere (copy format from row 2)
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A21:I21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("G21").Value = "esraa.sami@med.asu.edu.eg"
$ws.Range("H21").Value = "6/251"
$ws.Range("I21").Value = "Recorded"

# --- Row 28: G28 "Recorded By" list reordered ---
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"

# --- Row 30: G30 "Recorded By" list reordered ---
$ws.Range("G30").Value = "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
